# Cover-letter wording refresh:
#  - "managing ETL/Data warehouse & Business Intelligence products" ->
#    "working in ETL/DWH, Data Engineering, Data Analytics & Business
#     Intelligence products, managing"
#  - "responsible for product end to end development, data analytics,
#     testing, production implementation" ->
#    "responsible for end-to-end project implementation with requirement
#     gathering, estimation, development, testing, production implementation"
#  - "I am also experienced in startup organizations operations ... managed
#     product end to end lifecycle from requirement gathering to production
#     implementation along" ->
#    "I am experienced in HealthCare, ERP, Commodity Trade Management
#     domains, startup organizations operations ... managed product end to
#     end lifecycle along"
#  - "learn more about the position you have available" ->
#    "learn more about you have available"

$d = $word.ActiveDocument

$wdFindContinue = 1
$wdReplaceOne   = 1

$range1 = $d.Content
$range1.Find.Execute(
    "At my current place of employment, I am managing ETL/Data warehouse & Business Intelligence products, cross",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "At my current place of employment, I am working in ETL/DWH, Data Engineering, Data Analytics & Business Intelligence products, managing cross",
    $wdReplaceOne) | Out-Null

$range2 = $d.Content
$range2.Find.Execute(
    "responsible for product end to end development, data analytics, testing, production implementation",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "responsible for end-to-end project implementation with requirement gathering, estimation, development, testing, production implementation",
    $wdReplaceOne) | Out-Null

$range3 = $d.Content
$range3.Find.Execute(
    " I am also experienced in startup organizations operations and worked on IT infrastructure setup, developed ERP modules and managed product end to end lifecycle from requirement gathering to production implementation along",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    " I am experienced in HealthCare, ERP, Commodity Trade Management domains, startup organizations operations and worked on IT infrastructure setup, developed ERP modules and managed product end to end lifecycle along",
    $wdReplaceOne) | Out-Null

$range4 = $d.Content
$range4.Find.Execute(
    "Thank you for your time and consideration. I am eager to learn more about the position you have available",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "Thank you for your time and consideration. I am eager to learn more about you have available",
    $wdReplaceOne) | Out-Null
